$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 90 and 91, pushing the existing rows 90-179
# down to 92-181 (this also grows the used range to A1:T181, matching the
# new dimension in the target workbook).
$ws.Rows("90:91").Insert()

# The row that used to be "90" (Packham's Triumph / Primera, 07-May-2021)
# is now at row 92; the row that used to be "91" is now at row 93. Seed the
# two freshly inserted rows with copies of those so every untouched column
# (A,B,C,E,F,G,H,I,J,K,Q,R,T) carries the correct data/formatting, then
# overwrite just the cells that actually differ per the target data.
$ws.Range("A92:T92").Copy()
$ws.Range("A90:T90").PasteSpecial()

$ws.Range("A93:T93").Copy()
$ws.Range("A91:T91").PasteSpecial()

$excel.CutCopyMode = 0

# New row 90: 2022-04-07 (44658), Calidad "Especial", Volumen 50,
# Precio minimo/maximo/promedio 10000/10000/10000, Precio $/Kg 625.
$ws.Cells.Item(90, 4).Value = 44658
$ws.Cells.Item(90, 12).Value = "Especial"
$ws.Cells.Item(90, 13).Value = 50
$ws.Cells.Item(90, 14).Value = 10000
$ws.Cells.Item(90, 15).Value = 10000
$ws.Cells.Item(90, 16).Value = 10000
$ws.Cells.Item(90, 19).Value = 625

# New row 91: 2022-04-07 (44658), Volumen 80,
# Precio minimo/maximo/promedio 8000/9000/8500, Precio $/Kg 531.
$ws.Cells.Item(91, 4).Value = 44658
$ws.Cells.Item(91, 13).Value = 80
$ws.Cells.Item(91, 14).Value = 8000
$ws.Cells.Item(91, 15).Value = 9000
$ws.Cells.Item(91, 16).Value = 8500
$ws.Cells.Item(91, 19).Value = 531
